$d = $word.ActiveDocument

$replacements = @(
    @("92×95=", "76×22="),
    @("45×28=", "99×37="),
    @("11×60=", "27×30="),
    @("98×79=", "95×78="),
    @("45×38=", "62×86="),
    @("57×58=", "83×32="),
    @("40×38=", "11×72="),
    @("37×89=", "79×26="),
    @("54×49=", "68×86="),
    @("79×96=", "16×28="),
    @("80×66=", "59×80="),
    @("66×83=", "41×89="),
    @("11×46=", "73×62="),
    @("59×59=", "34×64="),
    @("61×24=", "38×56="),
    @("53×67=", "54×63="),
    @("35×48=", "35×61="),
    @("25×93=", "33×30="),
    @("44×94=", "77×68="),
    @("92×37=", "95×36="),
    @("14×78=", "25×66="),
    @("38×83=", "72×45="),
    @("63×51=", "42×87="),
    @("75×46=", "15×49="),
    @("93×55=", "93×44=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
